# Append a freshly-scraped batch (2026-01-28 12:44 JST) to the
# "ランサーズ" sheet:
#   - re-scored/re-ordered rows 2 & 3 pick up new title/price/score text
#   - rows 2-7 all pick up the new scrape timestamp in column A
#   - a brand-new row 8 is appended with its own hyperlink

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-28 12:44:14"

# --- Row 2: now the "大手SIer..." listing (score 375) ---------------------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5455098"
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- Row 3: now the "Dify..." listing (score 373, price re-quoted) --------
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "Difyと連携したAIチャットアプリ(Bubble)の実証実験用プロダクト(MVP)の開発"
$ws.Range("D3").Value = "80,000 円 ~ 90,000 円 / 募集期間 3 日、取引期間 0 日"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5480502"
$ws.Range("G3").Value = 373
$ws.Range("H3").Value = "🔥AI,Ai ◆開発 ◇アプリ"

# --- Rows 4-7: content unchanged, only the scrape timestamp refreshes -----
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp

# --- Row 8: brand-new listing appended -------------------------------------
$ws.Range("A8").Value = $newTimestamp
$ws.Range("B8").Value = "【継続】UTAGE実装者 _000ALL_RS"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "20,000 円 ~"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5480853"
$ws.Range("G8").Value = 10

# New hyperlink for F8, matching the style already used by F2:F7.
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5480853")
$ws.Range("F7").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5480853"
$excel.CutCopyMode = 0
